$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update non-price text columns (Coin name, Link, Volume label)
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("E43").Value = '42CEJICEJIBestin24h'

# Update Price column (D) values; force text format so the values are
# stored as strings (matching the original inline-string cell type)
# instead of being auto-converted to floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.38"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.15"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.277"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05799"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.501"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.132"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8165"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1362"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06938"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03124"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02873"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09404"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.733"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001515"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04685"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006004"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006267"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001236"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004624"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006891"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.505"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.138"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3192"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002329"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03668"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006252"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1056"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003396"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007474"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005260"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3697"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002213"
